# Generate LinkedIn draft + prompts
#
# Slide 1 ("NTPC Green Energy ... ") is rebuilt so that the picture is the
# first shape, followed by the headline textbox, followed by the body
# textbox -- and the headline/body boxes shrink a bit (auto-fit shrank once
# the picture moved up). Slides 2-5 just get the same headline/body resize.

# PowerPoint's Shape.Left/Top/Width/Height are single-precision (Single)
# COM properties; round-tripping an EMU value through pt = emu/12700 loses
# precision in the low bits, so nudge by a hair to land back on the exact
# EMU integer PowerPoint would have written.
function EmuToPt($emu) {
    return ($emu / 12700.0) + 0.00004
}

$p = $ppt.ActivePresentation

# --- Rebuild slide 1 -------------------------------------------------
# Insert a brand-new blank slide at position 1 (shapes added to it get a
# fresh id sequence: 2, 3, 4, ... in the order we add them -- matching the
# cNvPr ids in the target deck).
$newSlide = $p.Slides.Add(1, 7)
$oldSlide1 = $p.Slides.Item(2)

# 1) Picture first -> id 2 "Picture 1"
$oldSlide1.Shapes.Item(2).Copy()
$newPic = $newSlide.Shapes.Paste().Item(1)
$newPic.Name = "Picture 1"
$newPic.Left = EmuToPt(457200)
$newPic.Top = EmuToPt(1097280)
$newPic.Width = EmuToPt(8229600)
$newPic.Height = EmuToPt(5486400)

# 2) Headline textbox second -> id 3 "TextBox 2"
$oldSlide1.Shapes.Item(1).Copy()
$newTitle = $newSlide.Shapes.Paste().Item(1)
$newTitle.Name = "TextBox 2"
$newTitle.Left = EmuToPt(457200)
$newTitle.Top = EmuToPt(274320)
$newTitle.Width = EmuToPt(8229600)
$newTitle.Height = EmuToPt(731520)

# 3) Body textbox third -> id 4 "TextBox 3" (name unchanged)
$oldSlide1.Shapes.Item(3).Copy()
$newBody = $newSlide.Shapes.Paste().Item(1)
$newBody.Name = "TextBox 3"
$newBody.Left = EmuToPt(731520)
$newBody.Top = EmuToPt(4754880)
$newBody.Width = EmuToPt(7772400)
$newBody.Height = EmuToPt(1463040)

$oldSlide1.Delete()

# --- Slides 2-5: shrink headline + reposition/shrink body ------------
for ($i = 2; $i -le 5; $i++) {
    $s = $p.Slides.Item($i)

    $headline = $s.Shapes.Item(1)
    $headline.Height = EmuToPt(731520)

    $body = $s.Shapes.Item(2)
    $body.Top = EmuToPt(4754880)
    $body.Height = EmuToPt(1463040)
}
